$wb = $excel.ActiveWorkbook

$real = $wb.Worksheets.Item("RealDevices")

# Rename the two "real device" rows (columns A and B both hold the device/test-object name).
# Samsung (row 3) is written first so the new shared-string entries land in the same order
# as the target workbook.
$real.Range("A3").Value = "Samsung_Galaxy_S10_real"
$real.Range("B3").Value = "Samsung_Galaxy_S10_real"

$real.Range("A2").Value = "iPhone_XS_13_real"
$real.Range("B2").Value = "iPhone_XS_13_real"

# privateDevicesOnly becomes an explicit FALSE boolean (was blank) for the two renamed rows.
$real.Range("K2").Value = $false
$real.Range("K3").Value = $false

# "supported" column: rows 2 & 3 flip from true to false, row 5 flips from false to true.
# Use a leading apostrophe so Excel keeps these as text (shared-string) cells instead of
# auto-coercing "true"/"false" into native booleans.
$real.Range("L2").Value = "'false"
$real.Range("L3").Value = "'false"
$real.Range("L5").Value = "'true"

# Update the saved selection state on each sheet to match the edited workbook.
$browsers = $wb.Worksheets.Item("Browsers")
[void]$browsers.Activate()
[void]$browsers.Range("J2:J3").Select()

[void]$real.Activate()
[void]$real.Range("K9").Select()
